$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "May" budget row. Note: row 5 is intentionally left blank,
# the new data goes into row 6 (matching the source diff).
$ws.Range("A6").Value = "May"
$ws.Range("B6").Value = 7666
$ws.Range("B6").NumberFormat = $ws.Range("B4").NumberFormat

# Update the selected cell to mirror the author's final selection (B8)
$ws.Range("B8").Select()
